$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 224.14285
$ws.Range("I12").Value = 217.8
$ws.Range("K12").Value = 217.8
$ws.Range("M12").Value = -47.80000000000001
$ws.Range("H33").Value = 274
$ws.Range("I33").Value = 284.33334
$ws.Range("K33").Value = 284.33334
$ws.Range("M33").Value = -55.33334000000002
$ws.Range("H55").Value = 676.25
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H138").Value = 2373.31
$ws.Range("I138").Value = 2211.1304
$ws.Range("J138").Value = 2421.7532
$ws.Range("K138").Value = 6633.3912
$ws.Range("L138").Value = 7265.2596
$ws.Range("M138").Value = -1493.3912
$ws.Range("N138").Value = -17545.2596
$ws.Range("H141").Value = 4920.1
$ws.Range("I141").Value = 3600.1052
$ws.Range("K141").Value = 10800.3156
$ws.Range("M141").Value = -5620.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 5250
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5250
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5250
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -5538
$ws.Range("H32").Value = 3502.3062
$ws.Range("I32").Value = 3502.3062
$ws.Range("K32").Value = 3502.3062
$ws.Range("M32").Value = -3215.3062
$ws.Range("H61").Value = 3901.2173
$ws.Range("I61").Value = 2157.2727
$ws.Range("J61").Value = 5499.8335
$ws.Range("K61").Value = 2157.2727
$ws.Range("L61").Value = 5499.8335
$ws.Range("M61").Value = -1945.2727
$ws.Range("N61").Value = -5923.8335
$ws.Range("H74").Value = 243504
$ws.Range("I74").Value = 327981.2
$ws.Range("J74").Value = 4152
$ws.Range("K74").Value = 327981.2
$ws.Range("L74").Value = 4152
$ws.Range("M74").Value = -327107.2
$ws.Range("N74").Value = -5900
$ws.Range("H77").Value = 243504
$ws.Range("I77").Value = 327981.2
$ws.Range("J77").Value = 4152
$ws.Range("K77").Value = 1639906
$ws.Range("L77").Value = 20760
$ws.Range("M77").Value = -1635538
$ws.Range("N77").Value = -29496
$ws.Range("H88").Value = 4472.4443
$ws.Range("J88").Value = 6403.8
$ws.Range("L88").Value = 6403.8
$ws.Range("N88").Value = -7215.8
$ws.Range("H91").Value = 4472.4443
$ws.Range("J91").Value = 6403.8
$ws.Range("L91").Value = 6403.8
$ws.Range("N91").Value = -9211.799999999999
$ws.Range("H97").Value = 1024.3914
$ws.Range("I97").Value = 983.1
$ws.Range("K97").Value = 983.1
$ws.Range("M97").Value = -487.1
$ws.Range("H122").Value = 3316.024
$ws.Range("I122").Value = 3093.0527
$ws.Range("J122").Value = 5434.25
$ws.Range("K122").Value = 9279.158100000001
$ws.Range("L122").Value = 16302.75
$ws.Range("M122").Value = -6829.158100000001
$ws.Range("N122").Value = -21202.75
$ws.Range("H136").Value = 3901.2173
$ws.Range("I136").Value = 2157.2727
$ws.Range("J136").Value = 5499.8335
$ws.Range("K136").Value = 6471.8181
$ws.Range("L136").Value = 16499.5005
$ws.Range("M136").Value = -3921.8181
$ws.Range("N136").Value = -21599.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 71926.836
$ws.Range("I82").Value = 35248.5
$ws.Range("K82").Value = 35248.5
$ws.Range("M82").Value = -34865.5
$ws.Range("H85").Value = 71926.836
$ws.Range("I85").Value = 35248.5
$ws.Range("K85").Value = 35248.5
$ws.Range("M85").Value = -33922.5
$ws.Range("H94").Value = 117647930
$ws.Range("I94").Value = 125000300
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 125000300
$ws.Range("L94").Value = 10000
$ws.Range("M94").Value = -124999849
$ws.Range("N94").Value = -10902
$ws.Range("H97").Value = 10971.25
$ws.Range("I97").Value = 10971.25
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 10971.25
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -9980.25
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3769.535
$ws.Range("I31").Value = 2541.2727
$ws.Range("J31").Value = 7822.8
$ws.Range("K31").Value = 2541.2727
$ws.Range("L31").Value = 7822.8
$ws.Range("M31").Value = -2246.2727
$ws.Range("N31").Value = -8412.799999999999
$ws.Range("H34").Value = 3769.535
$ws.Range("I34").Value = 2541.2727
$ws.Range("J34").Value = 7822.8
$ws.Range("K34").Value = 2541.2727
$ws.Range("L34").Value = 7822.8
$ws.Range("M34").Value = -2339.2727
$ws.Range("N34").Value = -8226.799999999999
$ws.Range("H58").Value = 2886.862
$ws.Range("I58").Value = 1890.6875
$ws.Range("J58").Value = 4112.923
$ws.Range("K58").Value = 1890.6875
$ws.Range("L58").Value = 4112.923
$ws.Range("M58").Value = -1687.6875
$ws.Range("N58").Value = -4518.923
$ws.Range("H96").Value = 8722.666999999999
$ws.Range("J96").Value = 8722.666999999999
$ws.Range("L96").Value = 8722.666999999999
$ws.Range("N96").Value = -14214.667
$ws.Range("H122").Value = 3237.875
$ws.Range("I122").Value = 2022.2
$ws.Range("J122").Value = 5264
$ws.Range("K122").Value = 6066.6
$ws.Range("L122").Value = 15792
$ws.Range("M122").Value = -3616.6
$ws.Range("N122").Value = -20692
$ws.Range("H132").Value = 3582.12
$ws.Range("I132").Value = 3630.875
$ws.Range("J132").Value = 3495.4443
$ws.Range("K132").Value = 10892.625
$ws.Range("L132").Value = 10486.3329
$ws.Range("M132").Value = -8362.625
$ws.Range("N132").Value = -15546.3329
$ws.Range("H134").Value = 3611.0715
$ws.Range("I134").Value = 3396.261
$ws.Range("J134").Value = 4599.2
$ws.Range("K134").Value = 10188.783
$ws.Range("L134").Value = 13797.6
$ws.Range("M134").Value = -7653.782999999999
$ws.Range("N134").Value = -18867.6
$ws.Range("H136").Value = 2886.862
$ws.Range("I136").Value = 1890.6875
$ws.Range("J136").Value = 4112.923
$ws.Range("K136").Value = 5672.0625
$ws.Range("L136").Value = 12338.769
$ws.Range("M136").Value = -3122.0625
$ws.Range("N136").Value = -17438.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2893.625
$ws.Range("H92").Value = 5000
$ws.Range("J92").Value = 5000
$ws.Range("L92").Value = 15000
$ws.Range("N92").Value = -17496
$ws.Range("H137").Value = 2386.4
$ws.Range("I137").Value = 2138.2307
$ws.Range("J137").Value = 3999.5
$ws.Range("K137").Value = 6414.6921
$ws.Range("L137").Value = 11998.5
$ws.Range("M137").Value = -1314.6921
$ws.Range("N137").Value = -22198.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 711.5
$ws.Range("I2").Value = 964.5
$ws.Range("K2").Value = 964.5
$ws.Range("M2").Value = -851.5
$ws.Range("H102").Value = 1586.8
$ws.Range("I102").Value = 909.7143
$ws.Range("J102").Value = 3166.6667
$ws.Range("K102").Value = 909.7143
$ws.Range("L102").Value = 3166.6667
$ws.Range("M102").Value = 712.2857
$ws.Range("N102").Value = -6410.6667
$ws.Range("H113").Value = 2880.1333
$ws.Range("I113").Value = 2484.3076
$ws.Range("J113").Value = 5453
$ws.Range("K113").Value = 2484.3076
$ws.Range("L113").Value = 5453
$ws.Range("M113").Value = -314.3076000000001
$ws.Range("N113").Value = -9793
$ws.Range("H122").Value = 5079.857
$ws.Range("I122").Value = 4273.963
$ws.Range("J122").Value = 7799.75
$ws.Range("K122").Value = 12821.889
$ws.Range("L122").Value = 23399.25
$ws.Range("M122").Value = -10371.889
$ws.Range("N122").Value = -28299.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 901.46155
$ws.Range("I22").Value = 460
$ws.Range("K22").Value = 460
$ws.Range("M22").Value = -165
$ws.Range("H27").Value = 901.46155
$ws.Range("I27").Value = 460
$ws.Range("K27").Value = 460
$ws.Range("M27").Value = -353
$ws.Range("H40").Value = 4559.15
$ws.Range("I40").Value = 4457.237
$ws.Range("J40").Value = 6495.5
$ws.Range("K40").Value = 4457.237
$ws.Range("L40").Value = 6495.5
$ws.Range("M40").Value = -4321.237
$ws.Range("N40").Value = -6767.5
$ws.Range("H55").Value = 369.77777
$ws.Range("I55").Value = 224.66667
$ws.Range("K55").Value = 224.66667
$ws.Range("M55").Value = -51.66667000000001
$ws.Range("H82").Value = 899.2
$ws.Range("I82").Value = 899.2
$ws.Range("K82").Value = 899.2
$ws.Range("M82").Value = -538.2
$ws.Range("H85").Value = 899.2
$ws.Range("I85").Value = 899.2
$ws.Range("K85").Value = 899.2
$ws.Range("M85").Value = 348.8
$ws.Range("H95").Value = 35000
$ws.Range("J95").Value = 35000
$ws.Range("L95").Value = 35000
$ws.Range("N95").Value = -40492
$ws.Range("H99").Value = 24994.5
$ws.Range("I99").Value = 24994.5
$ws.Range("K99").Value = 24994.5
$ws.Range("M99").Value = -21999.5
$ws.Range("H132").Value = 5231.2573
$ws.Range("J132").Value = 10508.272
$ws.Range("L132").Value = 31524.816
$ws.Range("N132").Value = -36584.81600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 425.57693
$ws.Range("I113").Value = 435.33334
$ws.Range("K113").Value = 1306.00002
$ws.Range("M113").Value = 863.9999800000001
$ws.Range("H140").Value = 149000
$ws.Range("J140").Value = 149000
$ws.Range("L140").Value = 149000
$ws.Range("N140").Value = -159360
